$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = "O"
$ws.Range("H3").Value = "O"
$ws.Range("H4").Value = "O"
$ws.Range("H5").Value = "X"
$ws.Range("H6").Value = "O"
$ws.Range("H7").Value = "O"
$ws.Range("H8").Value = "O"
$ws.Range("H9").Value = "O"

$ws.Range("I5").Select()
